# Automatic update of files.
# Bump the "Förändrad" (Changed) date in column C, rows 2-43, from 45770 to 45771
# (i.e. from 2025-04-23 to 2025-04-24), leaving everything else untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 43; $row++) {
    $cell = $ws.Cells.Item($row, 3)   # Column C
    if ($cell.Value2 -eq 45770) {
        $cell.Value2 = 45771
    }
}
